$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.989.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.278.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.75"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +18.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.967"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +20.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.283.64"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "98.648.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000247"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.901.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.304.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +8.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.68"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000200"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.62"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.95"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.320"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +26.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.463.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.140"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +14.45%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.27"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +11.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.86"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.473"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.20%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.17"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "488.47"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.59"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.94"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "157.69"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.845"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.13%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.24"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +15.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.69"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.23%  "
